$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.520183
$ws.Range("H2").Value = 1.560549
$ws.Range("I2").Value = 0.03656880080220595
$ws.Range("J2").Value = 0.03656880080220595
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 8.601540095531998
$ws.Range("R2").Value = 77.41386085978799
$ws.Range("S2").Value = 0.00775424287244604
$ws.Range("T2").Value = 0.007754242872446042
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.520183
$ws.Range("H3").Value = 1.560549
$ws.Range("I3").Value = 0.03656880080220595
$ws.Range("J3").Value = 0.03656880080220595
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 21.13016152207867
$ws.Range("R3").Value = 190.171453698708
$ws.Range("S3").Value = 0.0190487287807124
$ws.Range("T3").Value = 0.0190487287807124
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.520183
$ws.Range("H4").Value = 1.560549
$ws.Range("I4").Value = 0.03656880080220595
$ws.Range("J4").Value = 0.03656880080220595
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 10.832930096907
$ws.Range("R4").Value = 97.496370872163
$ws.Range("S4").Value = 0.009765829149047511
$ws.Range("T4").Value = 0.009765829149047511
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.984906
$ws.Range("H5").Value = 38.954718
$ws.Range("I5").Value = 0.912837291778795
$ws.Range("J5").Value = 0.9128372917787949
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 214.713263593224
$ws.Range("R5").Value = 1932.419372339016
$ws.Range("S5").Value = 0.1935628707587173
$ws.Range("T5").Value = 0.1935628707587173
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.984906
$ws.Range("H6").Value = 38.954718
$ws.Range("I6").Value = 0.912837291778795
$ws.Range("J6").Value = 0.9128372917787949
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 527.4550708673841
$ws.Range("R6").Value = 4747.095637806457
$ws.Range("S6").Value = 0.475497954829445
$ws.Range("T6").Value = 0.475497954829445
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.984906
$ws.Range("H7").Value = 38.954718
$ws.Range("I7").Value = 0.912837291778795
$ws.Range("J7").Value = 0.9128372917787949
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 270.413640993474
$ws.Range("R7").Value = 2433.722768941266
$ws.Range("S7").Value = 0.2437764661906328
$ws.Range("T7").Value = 0.2437764661906327
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.719687
$ws.Range("H8").Value = 2.159061
$ws.Range("I8").Value = 0.05059390741899907
$ws.Range("J8").Value = 0.05059390741899907
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 11.900459235948
$ws.Range("R8").Value = 107.104133123532
$ws.Range("S8").Value = 0.01072820101799189
$ws.Range("T8").Value = 0.01072820101799189
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.719687
$ws.Range("H9").Value = 2.159061
$ws.Range("I9").Value = 0.05059390741899907
$ws.Range("J9").Value = 0.05059390741899907
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 29.23413982260134
$ws.Range("R9").Value = 263.107258403412
$ws.Range("S9").Value = 0.02635442232830478
$ws.Range("T9").Value = 0.02635442232830478
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.719687
$ws.Range("H10").Value = 2.159061
$ws.Range("I10").Value = 0.05059390741899907
$ws.Range("J10").Value = 0.05059390741899907
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 14.987646583323
$ws.Range("R10").Value = 134.888819249907
$ws.Range("S10").Value = 0.01351128407270241
$ws.Range("T10").Value = 0.01351128407270241
